$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7822.7676
$ws.Range("I62").Value = 7372.0557
$ws.Range("J62").Value = 10140.714
$ws.Range("K62").Value = 7372.0557
$ws.Range("L62").Value = 10140.714
$ws.Range("M62").Value = -6748.0557
$ws.Range("N62").Value = -11388.714
$ws.Range("H65").Value = 7822.7676
$ws.Range("I65").Value = 7372.0557
$ws.Range("J65").Value = 10140.714
$ws.Range("K65").Value = 36860.2785
$ws.Range("L65").Value = 50703.57
$ws.Range("M65").Value = -33740.2785
$ws.Range("N65").Value = -56943.57
$ws.Range("H76").Value = 4393.5557
$ws.Range("I76").Value = 3923.6667
$ws.Range("J76").Value = 5333.3335
$ws.Range("K76").Value = 3923.6667
$ws.Range("L76").Value = 5333.3335
$ws.Range("M76").Value = -3608.6667
$ws.Range("N76").Value = -5963.3335
$ws.Range("H79").Value = 4393.5557
$ws.Range("I79").Value = 3923.6667
$ws.Range("J79").Value = 5333.3335
$ws.Range("K79").Value = 3923.6667
$ws.Range("L79").Value = 5333.3335
$ws.Range("M79").Value = -2831.6667
$ws.Range("N79").Value = -7517.3335
$ws.Range("H106").Value = 2632
$ws.Range("I106").Value = 2507.6365
$ws.Range("K106").Value = 2507.6365
$ws.Range("M106").Value = -1876.6365
$ws.Range("H132").Value = 2245358.2
$ws.Range("J132").Value = 2163.3333
$ws.Range("L132").Value = 6489.999899999999
$ws.Range("N132").Value = -11549.9999
$ws.Range("H137").Value = 41550.5
$ws.Range("I137").Value = 46414.855
$ws.Range("K137").Value = 139244.565
$ws.Range("M137").Value = -136694.565
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17047.2
$ws.Range("I32").Value = 18038.275
$ws.Range("J32").Value = 5649.8335
$ws.Range("K32").Value = 18038.275
$ws.Range("L32").Value = 5649.8335
$ws.Range("M32").Value = -17751.275
$ws.Range("N32").Value = -6223.8335
$ws.Range("H61").Value = 9244.733
$ws.Range("I61").Value = 2344.75
$ws.Range("J61").Value = 17130.428
$ws.Range("K61").Value = 2344.75
$ws.Range("L61").Value = 17130.428
$ws.Range("M61").Value = -2132.75
$ws.Range("N61").Value = -17554.428
$ws.Range("H74").Value = 97640.22
$ws.Range("I74").Value = 102480.26
$ws.Range("K74").Value = 102480.26
$ws.Range("M74").Value = -101606.26
$ws.Range("H77").Value = 97640.22
$ws.Range("I77").Value = 102480.26
$ws.Range("K77").Value = 512401.3
$ws.Range("M77").Value = -508033.3
$ws.Range("H110").Value = 1978.5454
$ws.Range("I110").Value = 1983.9
$ws.Range("J110").Value = 1925
$ws.Range("K110").Value = 1983.9
$ws.Range("L110").Value = 1925
$ws.Range("M110").Value = 61.09999999999991
$ws.Range("N110").Value = -6015
$ws.Range("H131").Value = 78836
$ws.Range("J131").Value = 78836
$ws.Range("L131").Value = 78836
$ws.Range("N131").Value = -88916
$ws.Range("H132").Value = 1878.164
$ws.Range("I132").Value = 1677.826
$ws.Range("J132").Value = 2492.5334
$ws.Range("K132").Value = 5033.478
$ws.Range("L132").Value = 7477.600199999999
$ws.Range("M132").Value = -2503.478
$ws.Range("N132").Value = -12537.6002
$ws.Range("H136").Value = 9244.733
$ws.Range("I136").Value = 2344.75
$ws.Range("J136").Value = 17130.428
$ws.Range("K136").Value = 7034.25
$ws.Range("L136").Value = 51391.284
$ws.Range("M136").Value = -4484.25
$ws.Range("N136").Value = -56491.284
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 58938.215
$ws.Range("J20").Value = 1445.5555
$ws.Range("L20").Value = 1445.5555
$ws.Range("N20").Value = -1939.5555
$ws.Range("H105").Value = 4831.2085
$ws.Range("I105").Value = 5236.0557
$ws.Range("K105").Value = 5236.0557
$ws.Range("M105").Value = -3489.0557
$ws.Range("H107").Value = 41696.77
$ws.Range("I107").Value = 85083.5
$ws.Range("J107").Value = 4508.143
$ws.Range("K107").Value = 85083.5
$ws.Range("L107").Value = 4508.143
$ws.Range("M107").Value = -83163.5
$ws.Range("N107").Value = -8348.143
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5264981
$ws.Range("I31").Value = 5557202
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 5557202
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -5556907
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 5264981
$ws.Range("I34").Value = 5557202
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 5557202
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -5557000
$ws.Range("N34").Value = -5404
$ws.Range("H58").Value = 636
$ws.Range("I58").Value = 605.0833
$ws.Range("K58").Value = 605.0833
$ws.Range("M58").Value = -402.0833
$ws.Range("H86").Value = 44999.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 44999.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 44999.5
$ws.Range("N86").Value = -47245.5
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 44999.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 44999.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 224997.5
$ws.Range("N89").Value = -236229.5
$ws.Range("M89").ClearContents()
$ws.Range("H99").Value = 7771.5454
$ws.Range("I99").Value = 7525.25
$ws.Range("K99").Value = 7525.25
$ws.Range("M99").Value = -6027.25
$ws.Range("H122").Value = 1351.4166
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 7771.5454
$ws.Range("I126").Value = 7525.25
$ws.Range("K126").Value = 22575.75
$ws.Range("M126").Value = -20105.75
$ws.Range("H132").Value = 14174.195
$ws.Range("I132").Value = 17506.305
$ws.Range("J132").Value = 2178.6
$ws.Range("K132").Value = 52518.915
$ws.Range("L132").Value = 6535.799999999999
$ws.Range("M132").Value = -49988.915
$ws.Range("N132").Value = -11595.8
$ws.Range("H134").Value = 1337.8966
$ws.Range("I134").Value = 1242.8214
$ws.Range("K134").Value = 3728.4642
$ws.Range("M134").Value = -1193.4642
$ws.Range("H136").Value = 636
$ws.Range("I136").Value = 605.0833
$ws.Range("K136").Value = 1815.2499
$ws.Range("M136").Value = 734.7501
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1803
$ws.Range("I131").Value = 1361.125
$ws.Range("J131").Value = 2308
$ws.Range("K131").Value = 4083.375
$ws.Range("L131").Value = 6924
$ws.Range("M131").Value = 956.625
$ws.Range("N131").Value = -17004
$ws.Range("H132").Value = 1726.8462
$ws.Range("I132").Value = 2618.5
$ws.Range("J132").Value = 1459.35
$ws.Range("K132").Value = 23566.5
$ws.Range("L132").Value = 13134.15
$ws.Range("M132").Value = -21036.5
$ws.Range("N132").Value = -18194.15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3845.3057
$ws.Range("I70").Value = 4062.4285
$ws.Range("J70").Value = 3085.375
$ws.Range("K70").Value = 4062.4285
$ws.Range("L70").Value = 3085.375
$ws.Range("M70").Value = -3792.4285
$ws.Range("N70").Value = -3625.375
$ws.Range("H73").Value = 3845.3057
$ws.Range("I73").Value = 4062.4285
$ws.Range("J73").Value = 3085.375
$ws.Range("K73").Value = 4062.4285
$ws.Range("L73").Value = 3085.375
$ws.Range("M73").Value = -3126.4285
$ws.Range("N73").Value = -4957.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4824391.5
$ws.Range("I132").Value = 5574043
$ws.Range("J132").Value = 5201.4287
$ws.Range("K132").Value = 16722129
$ws.Range("L132").Value = 15604.2861
$ws.Range("M132").Value = -16719599
$ws.Range("N132").Value = -20664.2861
